# Auto-generated Excel COM-interop script applying numeric updates
# to the Adamantoise_Profits workbook (scheduled-runner price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce
$ws.Range("H11").Value = 27778744
$ws.Range("I11").Value = 27778744
$ws.Range("K11").Value = 27778744
$ws.Range("M11").Value = -27778604

# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 35721970
$ws.Range("I70").Value = 3400
$ws.Range("J70").Value = 50009396
$ws.Range("K70").Value = 10200
$ws.Range("L70").Value = 150028188
$ws.Range("M70").Value = -9930
$ws.Range("N70").Value = -150028728

# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 35721970
$ws.Range("I73").Value = 3400
$ws.Range("J73").Value = 50009396
$ws.Range("K73").Value = 10200
$ws.Range("L73").Value = 150028188
$ws.Range("M73").Value = -9264
$ws.Range("N73").Value = -150030060

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 1988.75
$ws.Range("I111").Value = 1988.75
$ws.Range("K111").Value = 5966.25
$ws.Range("M111").Value = -2899.25

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 3799.1765
$ws.Range("J112").Value = 4039.9333
$ws.Range("L112").Value = 12119.7999
$ws.Range("N112").Value = -14335.7999

# Row 116: Growing Up
$ws.Range("H116").Value = 56549.332
$ws.Range("I116").Value = 77324
$ws.Range("J116").Value = 15000
$ws.Range("K116").Value = 77324
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = -73882
$ws.Range("N116").Value = -21884

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3102.577
$ws.Range("I132").Value = 3140.8
$ws.Range("K132").Value = 9422.400000000001
$ws.Range("M132").Value = -6892.400000000001

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 4235.4707
$ws.Range("I137").Value = 2917.4167
$ws.Range("K137").Value = 8752.250100000001
$ws.Range("M137").Value = -6202.250100000001

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2613.7463
$ws.Range("I138").Value = 1740.7858
$ws.Range("J138").Value = 3240.487
$ws.Range("K138").Value = 5222.357400000001
$ws.Range("L138").Value = 9721.460999999999
$ws.Range("M138").Value = -82.35740000000078
$ws.Range("N138").Value = -20001.461

# Row 140: Tome for Tradition
$ws.Range("H140").Value = 199426.58
$ws.Range("I140").Value = 58994
$ws.Range("J140").Value = 222832
$ws.Range("K140").Value = 58994
$ws.Range("L140").Value = 222832
$ws.Range("M140").Value = -53814
$ws.Range("N140").Value = -233192

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 8016959
$ws.Range("I32").Value = 4387104.5
$ws.Range("J32").Value = 17869422
$ws.Range("K32").Value = 4387104.5
$ws.Range("L32").Value = 17869422
$ws.Range("M32").Value = -4386817.5
$ws.Range("N32").Value = -17869996

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2692.68
$ws.Range("I45").Value = 2274.65
$ws.Range("K45").Value = 2274.65
$ws.Range("M45").Value = -1897.65

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2953.7222
$ws.Range("I61").Value = 3068.9
$ws.Range("J61").Value = 2809.75
$ws.Range("K61").Value = 3068.9
$ws.Range("L61").Value = 2809.75
$ws.Range("M61").Value = -2856.9
$ws.Range("N61").Value = -3233.75

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 2462
$ws.Range("I74").Value = 2209.2856
$ws.Range("K74").Value = 2209.2856
$ws.Range("M74").Value = -1335.2856

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 2462
$ws.Range("I77").Value = 2209.2856
$ws.Range("K77").Value = 11046.428
$ws.Range("M77").Value = -6678.428

# Row 130: A Gift of Gloves
$ws.Range("H130").Value = 66001
$ws.Range("J130").Value = 66001
$ws.Range("L130").Value = 66001
$ws.Range("N130").Value = -76041

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2953.7222
$ws.Range("I136").Value = 3068.9
$ws.Range("J136").Value = 2809.75
$ws.Range("K136").Value = 9206.700000000001
$ws.Range("L136").Value = 8429.25
$ws.Range("M136").Value = -6656.700000000001
$ws.Range("N136").Value = -13529.25

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker
$ws.Range("H80").Value = 863.1111
$ws.Range("I80").Value = 800.3333
$ws.Range("J80").Value = 894.5
$ws.Range("K80").Value = 800.3333
$ws.Range("L80").Value = 894.5
$ws.Range("M80").Value = 197.6667
$ws.Range("N80").Value = -2890.5

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 863.1111
$ws.Range("I83").Value = 800.3333
$ws.Range("J83").Value = 894.5
$ws.Range("K83").Value = 4001.6665
$ws.Range("L83").Value = 4472.5
$ws.Range("M83").Value = 990.3334999999997
$ws.Range("N83").Value = -14456.5

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2589.25
$ws.Range("I105").Value = 1830.8889
$ws.Range("K105").Value = 1830.8889
$ws.Range("M105").Value = -83.88889999999992

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3603.2122
$ws.Range("I31").Value = 2011.1875
$ws.Range("J31").Value = 5101.5884
$ws.Range("K31").Value = 2011.1875
$ws.Range("L31").Value = 5101.5884
$ws.Range("M31").Value = -1716.1875
$ws.Range("N31").Value = -5691.5884

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3603.2122
$ws.Range("I34").Value = 2011.1875
$ws.Range("J34").Value = 5101.5884
$ws.Range("K34").Value = 2011.1875
$ws.Range("L34").Value = 5101.5884
$ws.Range("M34").Value = -1809.1875
$ws.Range("N34").Value = -5505.5884

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3615.3333
$ws.Range("I134").Value = 2428.1177
$ws.Range("K134").Value = 7284.353099999999
$ws.Range("M134").Value = -4749.353099999999

$ws = $wb.Worksheets.Item("CUL")
# Row 50: Moving Up in the World
$ws.Range("H50").Value = 1746.75
$ws.Range("I50").Value = 1134.9
$ws.Range("J50").Value = 2766.5
$ws.Range("K50").Value = 3404.7
$ws.Range("L50").Value = 8299.5
$ws.Range("M50").Value = -2923.7
$ws.Range("N50").Value = -9261.5

# Row 53: Rolanberry Fields Forever
$ws.Range("H53").Value = 1746.75
$ws.Range("I53").Value = 1134.9
$ws.Range("J53").Value = 2766.5
$ws.Range("K53").Value = 3404.7
$ws.Range("L53").Value = 8299.5
$ws.Range("M53").Value = -2923.7
$ws.Range("N53").Value = -9261.5

# Row 70: Persona non Gratin
$ws.Range("H70").Value = 4666.6665
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 24000
$ws.Range("M70").Value = -8685
$ws.Range("N70").Value = -24630

# Row 73: Recipe for Disaster (L)
$ws.Range("H73").Value = 4666.6665
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 24000
$ws.Range("M73").Value = -7908
$ws.Range("N73").Value = -26184

# Row 128: A Historical Flavor
$ws.Range("H128").Value = 210265.67
$ws.Range("I128").Value = 210265.67
$ws.Range("K128").Value = 630797.01
$ws.Range("M128").Value = -625817.01

# Row 132: More Mezcal
$ws.Range("H132").Value = 1920.0834
$ws.Range("J132").Value = 1920.0834
$ws.Range("L132").Value = 17280.7506
$ws.Range("N132").Value = -22340.7506

$ws = $wb.Worksheets.Item("GSM")
# Row 35: Necklet of Champions
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 1400.7885
$ws.Range("I102").Value = 1337.1277
$ws.Range("K102").Value = 1337.1277
$ws.Range("M102").Value = 284.8723

# Row 132: On Board for Lar
$ws.Range("H132").Value = 6099.6
$ws.Range("I132").Value = 5499.6665
$ws.Range("K132").Value = 16498.9995
$ws.Range("M132").Value = -13968.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 14046.549
$ws.Range("I7").Value = 13770.182
$ws.Range("J7").Value = 14722.111
$ws.Range("K7").Value = 13770.182
$ws.Range("L7").Value = 14722.111
$ws.Range("M7").Value = -13658.182
$ws.Range("N7").Value = -14946.111

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 338.45456
$ws.Range("I55").Value = 300.2353
$ws.Range("K55").Value = 300.2353
$ws.Range("M55").Value = -127.2353

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 1897.6666
$ws.Range("J61").Value = 1500
$ws.Range("L61").Value = 1500
$ws.Range("N61").Value = -1904

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 3891.25
$ws.Range("I68").Value = 3256.4285
$ws.Range("K68").Value = 3256.4285
$ws.Range("M68").Value = -2507.4285

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 3891.25
$ws.Range("I71").Value = 3256.4285
$ws.Range("K71").Value = 16282.1425
$ws.Range("M71").Value = -12538.1425

# Row 113: Peace in Rest
$ws.Range("H113").Value = 1897.6666
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

# Row 126: Battered Books
$ws.Range("H126").Value = 14046.549
$ws.Range("I126").Value = 13770.182
$ws.Range("J126").Value = 14722.111
$ws.Range("K126").Value = 41310.546
$ws.Range("L126").Value = 44166.333
$ws.Range("M126").Value = -38840.546
$ws.Range("N126").Value = -49106.333

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1136.5333
$ws.Range("I126").Value = 1153.4286
$ws.Range("K126").Value = 3460.2858
$ws.Range("M126").Value = -990.2857999999997

# Row 129: Lifetime of Gleaning
$ws.Range("H129").Value = 120000
$ws.Range("J129").Value = 120000
$ws.Range("L129").Value = 120000
$ws.Range("N129").Value = -130000

# Row 130: Skill Cap
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 131: A Better Bottom Line
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 4686.077
$ws.Range("I132").Value = 4004.625
$ws.Range("K132").Value = 12013.875
$ws.Range("M132").Value = -9483.875
